$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook lists test-case steps (TC2, TC3, TC4) in rows 18, 25, 32 (columns B & D).
# Target order after the edit (v1.0.2 -> v1.0.3):
#   TC2 (row 18): "analisar prestação de contas"
#   TC3 (row 25): "detalhar diária"
#   TC4 (row 32): "cancelar diária"

$ws.Range("B18").Value = "Beneficiário Clica em analisar prestação de contas."
$ws.Range("D18").Value = "SYSTEM Apresenta a tela de Analisar Prestação de Contas"

$ws.Range("B25").Value = "Beneficiário Clica em detalhar diária."
$ws.Range("D25").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"

$ws.Range("B32").Value = "Beneficiário Clica em cancelar diária."
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária"
